$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of rows 3 and 4 (G3, G4) from percentage values to "Done"
$ws.Range("G3").Value = "Done"
$ws.Range("G4").Value = "Done"

# Add new user story row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "user"
$ws.Range("C8").Value = "View the remaining hours of a student"
$ws.Range("D8").Value = "I can keep track of the student's remaining hours"
$ws.Range("F8").Value = "High"
$ws.Range("G8").Value = "Done"

# Update the selected cell / view to B9 (and clear the old topLeftCell/selection at G7)
$null = $ws.Range("B9").Select()
